$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Save" column - copy formatting (bold/border/alignment) from the
# neighboring "sum" header cell so it matches the other header cells' style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Cells.Item(1, 8).Value = "Save"

# Fill H2:H48 based on threshold of G column (rounded sum >= 10 => 1 else 0)
for ($r = 2; $r -le 48; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    if ([Math]::Round($g) -ge 10) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
